# Restructure the "Input" sheet to the new standard template column layout,
# clear the obsolete bold/centered header style, and drop the trailing
# empty "remarks" cells on the 갑지/을지 sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Helper: write a STRING value to a cell without letting Excel's COM layer
# auto-coerce date-shaped text (e.g. "2025-08-23") into a real date value.
# Force the cell to Text format first, assign, then strip the now-unneeded
# number format/style back off so the cell ends up with no "s" attribute,
# matching a freshly authored plain cell.
# NOTE: always pass a pre-computed variable as $val (never an inline
# "(...)" cast/expression at the call site) - this host's PowerShell layer
# silently drops inline cast expressions passed as call arguments and lets
# the original (uncast) COM variant flow through instead.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# --- Capture the old data before we overwrite anything ------------------
# NOTE: use .Value2 to read - .Value is broken on this host and returns a
# descriptor placeholder string instead of the real cell content.
$numRows = 5
$numOldCols = 17
$old = @{}
for ($r = 1; $r -le $numRows; $r++) {
    $old[$r] = @{}
    for ($c = 1; $c -le $numOldCols; $c++) {
        $old[$r][$c] = $ws.Cells.Item($r, $c).Value2
    }
}

# Old column index reference (1-based):
#  1 A 거래처명      2 B 현장명       3 C 발주일       4 D 납기일
#  5 E 발주번호      6 F 품목         7 G 규격         8 H 수량
#  9 I 단위         10 J 단가        11 K 공급가액    12 L 부가세
# 13 M 합계         14 N 대분류      15 O 중분류      16 P 소분류
# 17 Q 비고

# --- Clear the whole used range (values + formats) so the sheet shrinks
#     cleanly from 17 to 16 columns and loses the header cell styling ----
$ws.Range("A1:Q5").Clear()

# --- New header row (columns A-P), plain/no style ------------------------
$newHeaders = @(
    "발주일자", "납기일자", "거래처명", "거래처 이메일",
    "납품처명", "납품처 이메일", "프로젝트명", "대분류",
    "중분류", "소분류", "품목명", "규격",
    "수량", "단가", "총금액", "비고"
)
for ($c = 1; $c -le $newHeaders.Length; $c++) {
    $headerText = $newHeaders[$c - 1]
    $headerCell = $ws.Cells.Item(1, $c)
    Set-TextValue $headerCell $headerText
}

# --- New data rows, remapped from the captured old values ----------------
for ($r = 2; $r -le $numRows; $r++) {
    $vals = $old[$r]

    $companyName = $vals[1]            # old A 거래처명
    $siteName    = $vals[2]            # old B 현장명
    $orderDate   = $vals[3]            # old C 발주일
    $dueDate     = $vals[4]            # old D 납기일
    $bigCategory = $vals[14]           # old N 대분류
    $midCategory = $vals[15]           # old O 중분류
    $subCategory = $vals[16]           # old P 소분류
    $itemName    = $vals[6]            # old F 품목
    $spec        = $vals[7]            # old G 규격
    $qty         = $vals[8]            # old H 수량 (number)
    $unitPrice   = $vals[10]           # old J 단가 (number)
    $total       = $vals[13]           # old M 합계 (number)

    $companyEmail = "$companyName@example.com"
    $deliveryEmail = "delivery@example.com"

    $cellA = $ws.Cells.Item($r, 1)
    Set-TextValue $cellA $orderDate                    # 발주일자      <- old C 발주일
    $cellB = $ws.Cells.Item($r, 2)
    Set-TextValue $cellB $dueDate                      # 납기일자      <- old D 납기일
    $cellC = $ws.Cells.Item($r, 3)
    Set-TextValue $cellC $companyName                  # 거래처명      <- old A
    $cellD = $ws.Cells.Item($r, 4)
    Set-TextValue $cellD $companyEmail                 # 거래처 이메일 (new)
    $cellE = $ws.Cells.Item($r, 5)
    Set-TextValue $cellE $siteName                     # 납품처명      <- old B 현장명
    $cellF = $ws.Cells.Item($r, 6)
    Set-TextValue $cellF $deliveryEmail                # 납품처 이메일 (new)
    $cellG = $ws.Cells.Item($r, 7)
    Set-TextValue $cellG $siteName                     # 프로젝트명    <- old B 현장명
    $cellH = $ws.Cells.Item($r, 8)
    Set-TextValue $cellH $bigCategory                  # 대분류        <- old N
    $cellI = $ws.Cells.Item($r, 9)
    Set-TextValue $cellI $midCategory                  # 중분류        <- old O
    $cellJ = $ws.Cells.Item($r, 10)
    Set-TextValue $cellJ $subCategory                  # 소분류        <- old P
    $cellK = $ws.Cells.Item($r, 11)
    Set-TextValue $cellK $itemName                     # 품목명        <- old F 품목
    $cellL = $ws.Cells.Item($r, 12)
    Set-TextValue $cellL $spec                         # 규격          <- old G

    $ws.Cells.Item($r, 13).Value = $qty                # 수량          <- old H (number)
    $ws.Cells.Item($r, 14).Value = $unitPrice           # 단가          <- old J (number)
    $ws.Cells.Item($r, 15).Value = $total               # 총금액        <- old M 합계 (number)
    # 비고 (col 16) left blank
}

# --- 갑지 / 을지: drop the trailing empty "비고" cells (col I rows 2-5) ---
foreach ($sheetName in @("갑지", "을지")) {
    $s = $wb.Worksheets.Item($sheetName)
    for ($r = 2; $r -le 5; $r++) {
        $s.Cells.Item($r, 9).Clear()
    }
}
